$d = $word.ActiveDocument

function Find-ParagraphIndexByPrefix($prefix) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text.StartsWith($prefix)) {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------
# Change 1: paragraph "Each technical | request |  has:" (3 runs)
#           -> "Each technical request has:" (1 run, keeping the first
#           run's empty <w:rPr/> and its properties)
# ---------------------------------------------------------------------
$idx1 = Find-ParagraphIndexByPrefix "Each technical"
if ($idx1 -gt 0) {
    $p1 = $d.Paragraphs($idx1).Range
    $find1 = $p1.Duplicate
    $found1 = $find1.Find.Execute("Each technical ")
    if ($found1) {
        # Collapse to just after "Each technical " and extend through the
        # paragraph mark (End, not End-1 -- that is what makes the engine
        # fully re-coalesce the paragraph's run list instead of leaving the
        # old run boundaries in place), replacing "request has:" so the
        # trailing two runs merge into the leading run without touching
        # its (already empty) run properties.
        $mergeStart = $find1.End
        $mergeEnd = $p1.End
        $mergeRange = $d.Range($mergeStart, $mergeEnd)
        $mergeRange.Text = "request has:"
    }
}

# ---------------------------------------------------------------------
# Helper that turns the single word "service" inside a paragraph into
# "request", split across three runs: [prefix][request][suffix], all
# three ending up with an explicit (empty) <w:rPr/>, matching a real
# Word "retype with a transient formatting toggle" edit.
# ---------------------------------------------------------------------
function Split-ServiceToRequest($paragraphIndex) {
    $para = $d.Paragraphs($paragraphIndex).Range

    $findWord = $para.Duplicate
    $found = $findWord.Find.Execute("service")
    if (-not $found) { return }

    $wordStart = $findWord.Start
    $wordEnd = $findWord.End

    # Replace "service" with "request" in place first (same length).
    $wordRange = $d.Range($wordStart, $wordEnd)
    $wordRange.Text = "request"

    $paraStart = $para.Start
    $paraEnd = $para.End - 1   # exclude the paragraph mark

    # Toggling Bold on and back off on each of the three spans forces the
    # run to split from its neighbours while leaving an explicit empty
    # <w:rPr/> behind on all three resulting runs.
    $midRange = $d.Range($wordStart, $wordEnd)
    $midRange.Font.Bold = $true
    $midRange.Font.Bold = $false

    if ($wordStart -gt $paraStart) {
        $beforeRange = $d.Range($paraStart, $wordStart)
        $beforeRange.Font.Bold = $true
        $beforeRange.Font.Bold = $false
    }

    if ($wordEnd -lt $paraEnd) {
        $afterRange = $d.Range($wordEnd, $paraEnd)
        $afterRange.Font.Bold = $true
        $afterRange.Font.Bold = $false
    }
}

# ---------------------------------------------------------------------
# Change 2: "Blocks (...) one technical service applies (...)"
#           -> split into 3 runs around "request"
# ---------------------------------------------------------------------
$idx2 = Find-ParagraphIndexByPrefix "Blocks (one or more)"
if ($idx2 -gt 0) {
    Split-ServiceToRequest $idx2
}

# ---------------------------------------------------------------------
# Change 3: "Systems (...) one technical service applies (...)"
#           -> split into 3 runs around "request"
# ---------------------------------------------------------------------
$idx3 = Find-ParagraphIndexByPrefix "Systems (one or more)"
if ($idx3 -gt 0) {
    Split-ServiceToRequest $idx3
}
